$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (item 6, "BI ALCOFAN 150 MG 30 TABS."): balance ratio, sell price and
# transaction-count ratio were updated.
$ws.Range("H12").Value = "1:1"
$ws.Range("P12").Value = "296.4600"
$ws.Range("Q12").Value = "3:2"

# Row 23 (item 17, "MAALOX 20 ORAL SACHET SUSP."): same three fields updated.
$ws.Range("H23").Value = "0:12"
$ws.Range("P23").Value = "246.0000"
$ws.Range("Q23").Value = "2:1"

# Grand total at the bottom of the sell-price column reflects the two P-column
# changes above (+243.00 and +240.00 = +483.00).
$ws.Range("P44").Value = 4709.1899999999996

# Footer timestamp bumped by 2 minutes.
$ws.Range("A45").Value = "Wednesday, 24 September, 2025 7:57 PM"
